$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the K2 cell: it was a text "-" (shared string), should be a numeric 0
$ws.Range("K2").Value = 0

# Update the active selection to K3 (reflects where cursor was left on save)
$ws.Range("K3").Select()
